# Automatic map update (2025-09-08 07:30:22)
#
# The "NEW" sheet lists open work-order rows. This refresh drops one
# resolved/duplicate case (old row 57 - Caso 6002, LA PLATA AV. 832) and
# one case that's no longer tracked (old row 80 - Caso -582, Vilela 4019),
# which shifts every row in between up by one.
#
# Deleting entire rows (rather than rewriting every cell) lets Excel do
# the shifting for us and keeps every other row's values byte-for-byte
# identical.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 57 (Caso 6002) - everything below shifts up one row.
$ws.Rows.Item(57).Delete()

# After the shift above, the old row 80 (Caso -582, "Vilela 4019") is now
# at row 79 - remove it too, leaving the sheet with 78 data/header rows.
$ws.Rows.Item(79).Delete()
